$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the weekly price records for rows 25..66 down by one row (26..67),
# copying D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) from row (r-1) into row r.
# Work from the bottom up so we never overwrite a source row before reading it.
for ($r = 67; $r -ge 26; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value2  = $ws.Cells.Item($src, 4).Value2
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($src, 10).Value2
    $ws.Cells.Item($r, 11).Value2 = $ws.Cells.Item($src, 11).Value2
    $ws.Cells.Item($r, 12).Value2 = $ws.Cells.Item($src, 12).Value2
    $ws.Cells.Item($r, 13).Value2 = $ws.Cells.Item($src, 13).Value2
    $ws.Cells.Item($r, 16).Value2 = $ws.Cells.Item($src, 16).Value2

    # The rest of the row's columns (A,B,C,E,F,G,H,I,N,O,Q,R) are identical
    # across all these weekly records, so copy them down too to keep row 67
    # fully populated like the others (harmless no-op for rows 26..66).
    $ws.Cells.Item($r, 1).Value2  = $ws.Cells.Item($src, 1).Value2
    $ws.Cells.Item($r, 2).Value2  = $ws.Cells.Item($src, 2).Value2
    $ws.Cells.Item($r, 3).Value2  = $ws.Cells.Item($src, 3).Value2
    $ws.Cells.Item($r, 5).Value2  = $ws.Cells.Item($src, 5).Value2
    $ws.Cells.Item($r, 6).Value2  = $ws.Cells.Item($src, 6).Value2
    $ws.Cells.Item($r, 7).Value2  = $ws.Cells.Item($src, 7).Value2
    $ws.Cells.Item($r, 8).Value2  = $ws.Cells.Item($src, 8).Value2
    $ws.Cells.Item($r, 9).Value2  = $ws.Cells.Item($src, 9).Value2
    $ws.Cells.Item($r, 14).Value2 = $ws.Cells.Item($src, 14).Value2
    $ws.Cells.Item($r, 15).Value2 = $ws.Cells.Item($src, 15).Value2
    $ws.Cells.Item($r, 17).Value2 = $ws.Cells.Item($src, 17).Value2
    $ws.Cells.Item($r, 18).Value2 = $ws.Cells.Item($src, 18).Value2

    # Preserve the date-cell number format (column D) on the newly created row.
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat
}

# Row 25 becomes the newest weekly record.
$ws.Cells.Item(25, 4).Value2  = 44536
$ws.Cells.Item(25, 10).Value2 = 1000
$ws.Cells.Item(25, 11).Value2 = 3500
$ws.Cells.Item(25, 12).Value2 = 4000
$ws.Cells.Item(25, 13).Value2 = 3750
$ws.Cells.Item(25, 16).Value2 = 3750
